# Rename the 4th sheet ("3") to "Noble" and make it the active/selected tab
# (previously sheet "0" was the selected tab).

$wb = $excel.ActiveWorkbook

# Sheet "3" is the 4th sheet in the workbook (rId4 / sheetId 5) -> rename to "Noble"
$wsNoble = $wb.Worksheets.Item(4)
$wsNoble.Name = "Noble"

# Make the renamed "Noble" sheet the active tab (this also clears tabSelected
# on whichever sheet was previously active, i.e. sheet "0").
$wsNoble.Activate()
